# Corrected Calibration and Ingest Sheets for Coastal Gliders
# - Changed FLORT cal value for Scattering Angle (CC_scattering_angle) to 124
# - Changed FLORT cal value for Angular Resolution (CC_angular_resolution) to 1.076

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Row 4 holds CC_scattering_angle -> update value column (F) to 124
$ws.Range("F4").Value = 124

# Row 6 holds CC_angular_resolution -> update value column (F) to 1.076
$ws.Range("F6").Value = 1.076
